# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-370) from 2023-09-23 (serial 45192) to 2023-10-03 (serial 45202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C370").Value = 45202
